# Weekly update: a new price-report row for "Zapallo" (Camote) at Terminal
# Hortofrutícola Agro Chillán was inserted ahead of the existing series,
# pushing the previously-logged observations (old rows 64-88) down by one
# row (now rows 65-89). The sheet dimension grows from A1:R88 to A1:R89.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 64; this shifts the existing data
# (old rows 64..88) down to rows 65..89 and keeps column D's date
# formatting (style index 2) on the new row.
$ws.Rows("64:64").Insert()

# Populate the new row 64 with this period's observation. It follows the
# same Terminal Hortofrutícola Agro Chillán / Zapallo / Camote / "1a
# (guarda)" pattern as its neighbours, with its own date and prices.
$ws.Range("A64").Value = 7
$ws.Range("B64").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C64").Value = "Ñuble"
$ws.Range("D64").Value = 44455
$ws.Range("E64").Value = 16
$ws.Range("F64").Value = 100112045
$ws.Range("G64").Value = "Zapallo"
$ws.Range("H64").Value = "Camote"
$ws.Range("I64").Value = "1a (guarda)"
$ws.Range("J64").Value = 300
$ws.Range("K64").Value = 500
$ws.Range("L64").Value = 550
$ws.Range("M64").Value = 525
$ws.Range("N64").Value = "$/kilo (volumen en unidades)"
$ws.Range("O64").Value = "Región del Maule"
$ws.Range("P64").Value = 525
$ws.Range("Q64").Value = 1
$ws.Range("R64").Value = "Hortaliza"
